$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4333.3335
$ws.Range("I51").Value = 0
$ws.Range("J51").Value = 4333.3335
$ws.Range("K51").Value = 0
$ws.Range("L51").Value = 4333.3335
$ws.Range("N51").Value = -5301.3335
$ws.Range("M51").ClearContents()

$ws.Range("H88").Value = 2766.1667
$ws.Range("J88").Value = 4001.3333
$ws.Range("L88").Value = 4001.3333
$ws.Range("N88").Value = -4813.3333

$ws.Range("H91").Value = 2766.1667
$ws.Range("J91").Value = 4001.3333
$ws.Range("L91").Value = 4001.3333
$ws.Range("N91").Value = -6809.3333

$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = 0
$ws.Range("N123").ClearContents()

$ws.Range("H125").Value = 795.6667
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

$ws.Range("H134").Value = 80000
$ws.Range("J134").Value = 80000
$ws.Range("L134").Value = 80000
$ws.Range("N134").Value = -90140

$ws.Range("H138").Value = 6981.483
$ws.Range("I138").Value = 4517.2
$ws.Range("J138").Value = 7494.875
$ws.Range("K138").Value = 13551.6
$ws.Range("L138").Value = 22484.625
$ws.Range("M138").Value = -8411.599999999999
$ws.Range("N138").Value = -32764.625

$ws.Range("H140").Value = 90000
$ws.Range("J140").Value = 90000
$ws.Range("L140").Value = 90000
$ws.Range("N140").Value = -100360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1057.1666
$ws.Range("I2").Value = 968.7
$ws.Range("K2").Value = 968.7
$ws.Range("M2").Value = -855.7

$ws.Range("H61").Value = 2099.0908
$ws.Range("I61").Value = 1959
$ws.Range("K61").Value = 1959
$ws.Range("M61").Value = -1747

$ws.Range("H74").Value = 2168.6875
$ws.Range("I74").Value = 935.2222
$ws.Range("J74").Value = 3754.5715
$ws.Range("K74").Value = 935.2222
$ws.Range("L74").Value = 3754.5715
$ws.Range("M74").Value = -61.22220000000004
$ws.Range("N74").Value = -5502.5715

$ws.Range("H77").Value = 2168.6875
$ws.Range("I77").Value = 935.2222
$ws.Range("J77").Value = 3754.5715
$ws.Range("K77").Value = 4676.111
$ws.Range("L77").Value = 18772.8575
$ws.Range("M77").Value = -308.1109999999999
$ws.Range("N77").Value = -27508.8575

$ws.Range("H116").Value = 1057.1666
$ws.Range("I116").Value = 968.7
$ws.Range("K116").Value = 968.7
$ws.Range("M116").Value = 1325.3

$ws.Range("H132").Value = 1833.1136
$ws.Range("I132").Value = 1621.625
$ws.Range("K132").Value = 4864.875
$ws.Range("M132").Value = -2334.875

$ws.Range("H136").Value = 2099.0908
$ws.Range("I136").Value = 1959
$ws.Range("K136").Value = 5877
$ws.Range("M136").Value = -3327

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1057.1666
$ws.Range("I3").Value = 968.7
$ws.Range("K3").Value = 968.7
$ws.Range("M3").Value = -854.7

$ws.Range("H94").Value = 266.33334
$ws.Range("I94").Value = 129.5
$ws.Range("K94").Value = 129.5
$ws.Range("M94").Value = 321.5

$ws.Range("H97").Value = 8473.6
$ws.Range("I97").Value = 8473.6
$ws.Range("K97").Value = 8473.6
$ws.Range("M97").Value = -7482.6

$ws.Range("H134").Value = 858.9722
$ws.Range("I134").Value = 680.1177
$ws.Range("J134").Value = 3899.5
$ws.Range("K134").Value = 2040.3531
$ws.Range("L134").Value = 11698.5
$ws.Range("M134").Value = 494.6469
$ws.Range("N134").Value = -16768.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3935.5833
$ws.Range("I31").Value = 2287
$ws.Range("K31").Value = 2287
$ws.Range("M31").Value = -1992

$ws.Range("H34").Value = 3935.5833
$ws.Range("I34").Value = 2287
$ws.Range("K34").Value = 2287
$ws.Range("M34").Value = -2085

$ws.Range("H99").Value = 11077.28
$ws.Range("I99").Value = 8782.076999999999
$ws.Range("J99").Value = 13563.75
$ws.Range("K99").Value = 8782.076999999999
$ws.Range("L99").Value = 13563.75
$ws.Range("M99").Value = -7284.076999999999
$ws.Range("N99").Value = -16559.75

$ws.Range("H126").Value = 11077.28
$ws.Range("I126").Value = 8782.076999999999
$ws.Range("J126").Value = 13563.75
$ws.Range("K126").Value = 26346.231
$ws.Range("L126").Value = 40691.25
$ws.Range("M126").Value = -23876.231
$ws.Range("N126").Value = -45631.25

$ws.Range("H134").Value = 1509.0465
$ws.Range("I134").Value = 1159.9412
$ws.Range("J134").Value = 2827.889
$ws.Range("K134").Value = 3479.8236
$ws.Range("L134").Value = 8483.667000000001
$ws.Range("M134").Value = -944.8235999999997
$ws.Range("N134").Value = -13553.667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1617.3334
$ws.Range("J34").Value = 1974.4
$ws.Range("L34").Value = 5923.200000000001
$ws.Range("N34").Value = -6091.200000000001

$ws.Range("H60").Value = 765
$ws.Range("I60").Value = 765
$ws.Range("K60").Value = 2295
$ws.Range("M60").Value = -2044

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 7999.5
$ws.Range("I55").Value = 7999
$ws.Range("K55").Value = 7999
$ws.Range("M55").Value = -7672

$ws.Range("H113").Value = 4714.143
$ws.Range("I113").Value = 2999
$ws.Range("K113").Value = 2999
$ws.Range("M113").Value = -829

$ws.Range("H132").Value = 2607.25
$ws.Range("I132").Value = 2116.4
$ws.Range("J132").Value = 4079.8
$ws.Range("K132").Value = 6349.200000000001
$ws.Range("L132").Value = 12239.4
$ws.Range("M132").Value = -3819.200000000001
$ws.Range("N132").Value = -17299.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 6622.8335
$ws.Range("J61").Value = 4001
$ws.Range("L61").Value = 4001
$ws.Range("N61").Value = -4405

$ws.Range("H113").Value = 6622.8335
$ws.Range("J113").Value = 4001
$ws.Range("L113").Value = 4001
$ws.Range("N113").Value = -8341

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1320.7142
$ws.Range("I107").Value = 631.7778
$ws.Range("J107").Value = 2560.8
$ws.Range("K107").Value = 1895.3334
$ws.Range("L107").Value = 7682.400000000001
$ws.Range("M107").Value = 24.66660000000002
$ws.Range("N107").Value = -11522.4

$ws.Range("H113").Value = 2535.7222
$ws.Range("I113").Value = 2174.5
$ws.Range("J113").Value = 2716.3333
$ws.Range("K113").Value = 6523.5
$ws.Range("L113").Value = 8148.999899999999
$ws.Range("M113").Value = -4353.5
$ws.Range("N113").Value = -12488.9999

$ws.Range("H132").Value = 2231.08
$ws.Range("I132").Value = 1830.4736
$ws.Range("J132").Value = 3499.6667
$ws.Range("K132").Value = 5491.4208
$ws.Range("L132").Value = 10499.0001
$ws.Range("M132").Value = -2961.4208
$ws.Range("N132").Value = -15559.0001
